$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------
# 1. Reorganize columns.
#    Before: A=EN, B=pics, C=ES, D=FR, E=sizeW, F=sizeH
#    After : A=pics, B=EN, C=ES, D=FR, E=DE(new), F=sizeW, G=sizeH
# --------------------------------------------------------------------

# Swap A and B so the "pics" column becomes the first column.
$ws.Columns.Item(2).Cut() | Out-Null
$ws.Columns.Item(1).Insert() | Out-Null

# Insert a brand new column before the (now) sizeW column to host the
# German translation.
$ws.Columns.Item(5).Insert() | Out-Null

# --------------------------------------------------------------------
# 2. Header row renames (row 1) - keys were renamed.
# --------------------------------------------------------------------
$ws.Range("A1").Value = "inst_pics"
$ws.Range("B1").Value = "inst_msg_EN"
$ws.Range("C1").Value = "inst_msg_ES"
$ws.Range("D1").Value = "inst_msg_FR"
$ws.Range("E1").Value = "inst_msg_DE"
$ws.Range("F1").Value = "image_w"
$ws.Range("G1").Value = "image_h"

# --------------------------------------------------------------------
# 3. Fill in the new German column (rows 2-4).
# --------------------------------------------------------------------
$ws.Range("E2").Value = "Anweisungen:`r`nIn dieser Aufgabe sehen Sie Bilder von linken oder rechten Händen, deren Handflächen nach oben oder unten zeigen. Die Bilder werden in verschiedenen Winkeln gedreht.`r`nIhre Aufgabe ist es zu bestimmen, ob das Bild einer linken oder rechten Hand entspricht.`r`nIhr Ziel ist es, sowohl SCHNELL als auch GENAU zu antworten.`r`nJedes Bild wird angezeigt, bis Sie geantwortet haben. Das nächste Bild erscheint automatisch."

$ws.Range("E3").Value = "Bitte verwenden Sie nur Ihren ZEIGEFINGER und MITTELFINGER der RECHTEN HAND, um zu antworten.`r`nLegen Sie Ihren Zeigefinger auf die Taste „G“ und den Mittelfinger auf die Taste „H“ Ihrer Tastatur.`r`nZum Antworten:`r`nLinke Hand = G | H = Rechte Hand`r`nSie müssen Ihre Hand während der gesamten Aufgabe auf der Tastatur halten.`r`nHalten Sie Ihre andere Hand auf dem Tisch, in derselben Position und so ruhig wie möglich."

$ws.Range("E4").Value = "Nach jedem Bild erhalten Sie ein kurzes Feedback zu Ihrer Antwort:`r`nWenn Sie korrekt antworten, wird das entsprechende Feld grün.`r`nWenn Sie falsch antworten, wird das entsprechende Feld rot.`r`nDenken Sie daran, dass Ihr Ziel darin besteht, so genau und schnell wie möglich zu antworten."

# --------------------------------------------------------------------
# 4. Style the header row: solid blue fill + white font.
#    Build the style on a scratch cell first so only a single new
#    cellXf (font+fill) is produced, then copy it onto the header row
#    with a formats-only paste (keeps text / values untouched).
# --------------------------------------------------------------------
$scratch = $ws.Range("J1")
$scratch.Interior.Color = 15773696
$scratch.Font.ThemeColor = 2
$scratch.Copy() | Out-Null

$header = $ws.Range("A1:G1")
$header.PasteSpecial(-4122) | Out-Null

$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false

# --------------------------------------------------------------------
# 5. Column widths for the final layout (A & B already match after the
#    swap above, so only C, D, E, F, G need adjusting).
# --------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 13.916666666666668
$ws.Columns.Item(4).ColumnWidth = 14.583333333333332
$ws.Columns.Item(5).ColumnWidth = 14.583333333333332
$ws.Columns.Item(6).ColumnWidth = 7.916666666666666
$ws.Columns.Item(7).ColumnWidth = 7.583333333333334

# --------------------------------------------------------------------
# 6. Keep the original row heights (wrapping the new German text would
#    otherwise trigger auto row-height growth).
# --------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 24.5
$ws.Rows.Item(3).RowHeight = 24.5
$ws.Rows.Item(4).RowHeight = 24.5

# --------------------------------------------------------------------
# 7. Selection matches the saved view in the target file.
# --------------------------------------------------------------------
$ws.Range("E4").Select() | Out-Null
